$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 56; this shifts all existing rows (56..87) down by one
# (to 57..88) and carries their formatting/values with them.
$ws.Rows.Item(56).Insert()

# Copy the date-cell format (style) from the row above into the new D56 cell,
# since Excel's Insert normally pulls formatting from the row above anyway,
# but we set it explicitly to be safe.
$ws.Range("D55").Copy() | Out-Null
$ws.Range("D56").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row 56 with the new weekly data point.
$ws.Cells.Item(56, 1).Value = 11
$ws.Cells.Item(56, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(56, 3).Value = "Bíobío"
$ws.Cells.Item(56, 4).Value = 44572
$ws.Cells.Item(56, 5).Value = 8
$ws.Cells.Item(56, 6).Value = "Fruta"
$ws.Cells.Item(56, 7).Value = 100108
$ws.Cells.Item(56, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(56, 9).Value = 100108002
$ws.Cells.Item(56, 10).Value = "Mango"
$ws.Cells.Item(56, 11).Value = "Sin especificar"
$ws.Cells.Item(56, 12).Value = "Primera"
$ws.Cells.Item(56, 13).Value = 200
$ws.Cells.Item(56, 14).Value = 6000
$ws.Cells.Item(56, 15).Value = 6500
$ws.Cells.Item(56, 16).Value = 6250
$ws.Cells.Item(56, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(56, 18).Value = "Perú"
$ws.Cells.Item(56, 19).Value = 1562
$ws.Cells.Item(56, 20).Value = 4
